$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2398107.8
$ws.Range("J9").Value = 7993265.5
$ws.Range("L9").Value = 7993265.5
$ws.Range("N9").Value = -7993603.5

# Sheet ALC, Row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3327.9092
$ws.Range("I40").Value = 2734.3333
$ws.Range("J40").Value = 4040.2
$ws.Range("K40").Value = 2734.3333
$ws.Range("L40").Value = 4040.2
$ws.Range("M40").Value = -2559.3333
$ws.Range("N40").Value = -4390.2

# Sheet ALC, Row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2328310
$ws.Range("I70").Value = 2749.5
$ws.Range("K70").Value = 8248.5
$ws.Range("M70").Value = -7978.5

# Sheet ALC, Row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2328310
$ws.Range("I73").Value = 2749.5
$ws.Range("K73").Value = 8248.5
$ws.Range("M73").Value = -7312.5

# Sheet ALC, Row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3932.6667
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 3932.6667
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 11798.0001
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -21878.0001

# Sheet ALC, Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 997.7353000000001
$ws.Range("I132").Value = 652.65515
$ws.Range("K132").Value = 1957.96545
$ws.Range("M132").Value = 572.0345499999999

# Sheet ARM, Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23220.896
$ws.Range("I32").Value = 23089.129
$ws.Range("K32").Value = 23089.129
$ws.Range("M32").Value = -22802.129

# Sheet ARM, Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1909.3125
$ws.Range("I45").Value = 1605.3
$ws.Range("J45").Value = 2416
$ws.Range("K45").Value = 1605.3
$ws.Range("L45").Value = 2416
$ws.Range("M45").Value = -1228.3
$ws.Range("N45").Value = -3170

# Sheet ARM, Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6259.0356
$ws.Range("I61").Value = 5018.364
$ws.Range("J61").Value = 10808.167
$ws.Range("K61").Value = 5018.364
$ws.Range("L61").Value = 10808.167
$ws.Range("M61").Value = -4806.364
$ws.Range("N61").Value = -11232.167

# Sheet ARM, Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5399.4
$ws.Range("I122").Value = 4061.75
$ws.Range("K122").Value = 12185.25
$ws.Range("M122").Value = -9735.25

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3720.073
$ws.Range("I132").Value = 2633.2056
$ws.Range("J132").Value = 7169.696
$ws.Range("K132").Value = 7899.6168
$ws.Range("L132").Value = 21509.088
$ws.Range("M132").Value = -5369.6168
$ws.Range("N132").Value = -26569.088

# Sheet ARM, Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6259.0356
$ws.Range("I136").Value = 5018.364
$ws.Range("J136").Value = 10808.167
$ws.Range("K136").Value = 15055.092
$ws.Range("L136").Value = 32424.501
$ws.Range("M136").Value = -12505.092
$ws.Range("N136").Value = -37524.501

# Sheet CRP, Row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6377.467
$ws.Range("I134").Value = 4440
$ws.Range("J134").Value = 8591.714
$ws.Range("K134").Value = 13320
$ws.Range("L134").Value = 25775.142
$ws.Range("M134").Value = -10785
$ws.Range("N134").Value = -30845.142

# Sheet CUL, Row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1218.579
$ws.Range("I55").Value = 978
$ws.Range("J55").Value = 1304.5
$ws.Range("K55").Value = 2934
$ws.Range("L55").Value = 3913.5
$ws.Range("M55").Value = -2757
$ws.Range("N55").Value = -4267.5

# Sheet CUL, Row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5858
$ws.Range("J63").Value = 11500
$ws.Range("L63").Value = 34500
$ws.Range("N63").Value = -35998

# Sheet CUL, Row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 5858
$ws.Range("J66").Value = 11500
$ws.Range("L66").Value = 103500
$ws.Range("N66").Value = -110988

# Sheet CUL, Row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 316
$ws.Range("I75").Value = 83
$ws.Range("J75").Value = 1015
$ws.Range("K75").Value = 249
$ws.Range("L75").Value = 3045
$ws.Range("M75").Value = 749
$ws.Range("N75").Value = -5041

# Sheet CUL, Row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 316
$ws.Range("I78").Value = 83
$ws.Range("J78").Value = 1015
$ws.Range("K78").Value = 747
$ws.Range("L78").Value = 9135
$ws.Range("M78").Value = 4245
$ws.Range("N78").Value = -19119

# Sheet CUL, Row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 4676.375
$ws.Range("I87").Value = 4201.7144
$ws.Range("K87").Value = 12605.1432
$ws.Range("M87").Value = -11357.1432

# Sheet CUL, Row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 4676.375
$ws.Range("I90").Value = 4201.7144
$ws.Range("K90").Value = 37815.4296
$ws.Range("M90").Value = -31575.4296

# Sheet CUL, Row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1297.5518
$ws.Range("I107").Value = 473.76923
$ws.Range("J107").Value = 1966.875
$ws.Range("K107").Value = 1421.30769
$ws.Range("L107").Value = 5900.625
$ws.Range("M107").Value = 498.6923099999999
$ws.Range("N107").Value = -9740.625

# Sheet CUL, Row 120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 1029.5
$ws.Range("I120").Value = 1029.5
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 3088.5
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = 1749.5
$ws.Range("N120").ClearContents()

# Sheet GSM, Row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622

# Sheet GSM, Row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12403
$ws.Range("I122").Value = 7908
$ws.Range("K122").Value = 23724
$ws.Range("M122").Value = -21274

# Sheet LTW, Row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2523
$ws.Range("I7").Value = 1984.4
$ws.Range("J7").Value = 3196.25
$ws.Range("K7").Value = 1984.4
$ws.Range("L7").Value = 3196.25
$ws.Range("M7").Value = -1872.4
$ws.Range("N7").Value = -3420.25

# Sheet LTW, Row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6818.3706
$ws.Range("J46").Value = 7927.273
$ws.Range("L46").Value = 7927.273
$ws.Range("N46").Value = -8303.273000000001

# Sheet LTW, Row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3646.1875
$ws.Range("I68").Value = 2193.625
$ws.Range("J68").Value = 5098.75
$ws.Range("K68").Value = 2193.625
$ws.Range("L68").Value = 5098.75
$ws.Range("M68").Value = -1444.625
$ws.Range("N68").Value = -6596.75

# Sheet LTW, Row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3646.1875
$ws.Range("I71").Value = 2193.625
$ws.Range("J71").Value = 5098.75
$ws.Range("K71").Value = 10968.125
$ws.Range("L71").Value = 25493.75
$ws.Range("M71").Value = -7224.125
$ws.Range("N71").Value = -32981.75

# Sheet LTW, Row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2523
$ws.Range("I126").Value = 1984.4
$ws.Range("J126").Value = 3196.25
$ws.Range("K126").Value = 5953.200000000001
$ws.Range("L126").Value = 9588.75
$ws.Range("M126").Value = -3483.200000000001
$ws.Range("N126").Value = -14528.75

# Sheet WVR, Row 24
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 40000
$ws.Range("I24").Value = 40000
$ws.Range("K24").Value = 40000
$ws.Range("M24").Value = -39770

# Sheet WVR, Row 28
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 45000
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Sheet WVR, Row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1057.12
$ws.Range("I100").Value = 1093.3478
$ws.Range("J100").Value = 640.5
$ws.Range("K100").Value = 2186.6956
$ws.Range("L100").Value = 1281
$ws.Range("M100").Value = -1645.6956
$ws.Range("N100").Value = -2363

# Sheet WVR, Row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2172.4856
$ws.Range("I122").Value = 1905.8462
$ws.Range("K122").Value = 5717.5386
$ws.Range("M122").Value = -3267.5386

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3941.7778
$ws.Range("I132").Value = 2592.739
$ws.Range("K132").Value = 7778.217000000001
$ws.Range("M132").Value = -5248.217000000001

Write-Output "Applied all changes"